$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new headers in row 1, matching the style of the existing headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), one pair per row from row 2 to row 64.
$data = @(
    @(7,7),@(1,1),@(8,8),@(8,8),@(6,6),@(6,7),@(9,10),@(7,7),@(7,8),@(11,11),
    @(6,7),@(8,8),@(8,8),@(7,7),@(7,7),@(6,6),@(8,8),@(9,9),@(8,8),@(7,7),
    @(8,8),@(8,8),@(8,8),@(7,7),@(6,7),@(8,8),@(7,7),@(7,8),@(10,10),@(7,7),
    @(6,6),@(5,6),@(7,8),@(1,1),@(6,7),@(4,5),@(1,2),@(1,1),@(6,6),@(5,6),
    @(6,6),@(1,2),@(1,1),@(2,4),@(6,6),@(7,7),@(5,6),@(5,6),@(7,8),@(7,7),
    @(4,4),@(10,11),@(7,7),@(6,7),@(7,8),@(8,9),@(6,6),@(9,9),@(5,6),@(6,6),
    @(6,6),@(5,5),@(4,4)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
